$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G header
$ws.Range("G1").Value = "Última fecha"

# Row -> last-date value ("Última fecha") for column G
$dates = @{
    2  = "2024-02-16 00:00:00"
    3  = "2024-02-16 00:00:00"
    4  = "2024-01-26 00:00:00"
    5  = "2024-02-16 00:00:00"
    6  = "2024-02-16 00:00:00"
    7  = "2024-02-16 00:00:00"
    8  = "2024-01-26 00:00:00"
    9  = "2024-02-16 00:00:00"
    10 = "2024-01-26 00:00:00"
    11 = "2024-01-26 00:00:00"
    12 = "2024-01-26 00:00:00"
    13 = "2024-01-26 00:00:00"
    14 = "2024-01-26 00:00:00"
    15 = "2024-01-26 00:00:00"
    16 = "2024-02-16 00:00:00"
    17 = "2024-01-26 00:00:00"
    18 = "2024-02-16 00:00:00"
    19 = "2024-01-26 00:00:00"
    20 = "2023-12-06 00:00:00"
    21 = "2024-01-26 00:00:00"
    22 = "2024-01-26 00:00:00"
    23 = "2024-01-26 00:00:00"
    24 = "2024-01-26 00:00:00"
    25 = "2024-02-19 00:00:00"
    26 = "2024-01-26 00:00:00"
    27 = "2024-01-26 00:00:00"
    28 = "2024-02-16 00:00:00"
    29 = "2024-01-26 00:00:00"
    30 = "2024-01-26 00:00:00"
    31 = "2024-01-26 00:00:00"
    32 = "2024-01-26 00:00:00"
    33 = "2024-01-29 00:00:00"
    34 = "2024-01-26 00:00:00"
    35 = "2024-01-26 00:00:00"
    36 = "2024-01-26 00:00:00"
    37 = "2024-01-26 00:00:00"
    38 = "2024-01-26 00:00:00"
    39 = "2024-01-26 00:00:00"
    40 = "2024-01-11 00:00:00"
    41 = "2024-02-07 00:00:00"
    42 = "2024-01-04 00:00:00"
    43 = "2024-01-26 00:00:00"
    44 = "2024-01-26 00:00:00"
    45 = "2024-01-26 00:00:00"
    46 = "2024-01-26 00:00:00"
}

foreach ($row in $dates.Keys) {
    $ws.Cells.Item($row, 7).Value = $dates[$row]
}

# Column F: every cell that still carries a numeric price becomes "Sin precio"
for ($row = 2; $row -le 46; $row++) {
    $ws.Cells.Item($row, 6).Value = "Sin precio"
}
